$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030359391618642
$ws.Range("D2").Value = 1.032743152920536
$ws.Range("E2").Value = 1.039297547903983
$ws.Range("F2").Value = 1.048419981340831
$ws.Range("I2").Value = 1.032490467137752
$ws.Range("J2").Value = 1.035501071422004
$ws.Range("K2").Value = 1.035547504045563
$ws.Range("L2").Value = 1.04208313535606
$ws.Range("M2").Value = 1.05117987703966
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031302376646201
$ws.Range("D3").Value = 1.033406807290033
$ws.Range("E3").Value = 1.040142136546105
$ws.Range("F3").Value = 1.049364256358673
$ws.Range("I3").Value = 1.032627378061937
$ws.Range("J3").Value = 1.03608551039711
$ws.Range("K3").Value = 1.036020382397979
$ws.Range("L3").Value = 1.042737823025312
$ws.Range("M3").Value = 1.051935851555375
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031913107913507
$ws.Range("D4").Value = 1.033836476435535
$ws.Range("E4").Value = 1.04068951334721
$ws.Range("F4").Value = 1.049976232178579
$ws.Range("I4").Value = 1.032714746184533
$ws.Range("J4").Value = 1.036463630344673
$ws.Range("K4").Value = 1.036325924803669
$ws.Range("L4").Value = 1.043161676308438
$ws.Range("M4").Value = 1.052425358751945
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03216999158909
$ws.Range("D5").Value = 1.03401716539184
$ws.Range("E5").Value = 1.040919837749432
$ws.Range("F5").Value = 1.050233736777928
$ws.Range("I5").Value = 1.03275118253241
$ws.Range("J5").Value = 1.036622578676288
$ws.Range("K5").Value = 1.036454268057021
$ws.Range("L5").Value = 1.043339917317602
$ws.Range("M5").Value = 1.052631227955123
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032213131212758
$ws.Range("D6").Value = 1.03404750710561
$ws.Range("E6").Value = 1.040958522338286
$ws.Range("F6").Value = 1.050276986375632
$ws.Range("I6").Value = 1.032757283145095
$ws.Range("J6").Value = 1.036649265990513
$ws.Range("K6").Value = 1.03647581116144
$ws.Range("L6").Value = 1.043369847846927
$ws.Range("M6").Value = 1.052665798976111
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031916539887431
$ws.Range("D7").Value = 1.033838890591875
$ws.Range("E7").Value = 1.040692590140813
$ws.Range("F7").Value = 1.049979672065882
$ws.Range("I7").Value = 1.03271523420162
$ws.Range("J7").Value = 1.036465754273612
$ws.Range("K7").Value = 1.036327640153402
$ws.Range("L7").Value = 1.043164057766372
$ws.Range("M7").Value = 1.052428109272102
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.03067796186836
$ws.Range("D8").Value = 1.03296738734446
$ws.Range("E8").Value = 1.039582799281087
$ws.Range("F8").Value = 1.04873890244646
$ws.Range("I8").Value = 1.03253698959124
$ws.Range("J8").Value = 1.035698595269883
$ws.Range("K8").Value = 1.035707406269241
$ws.Range("L8").Value = 1.042304342508497
$ws.Range("M8").Value = 1.05143529075626
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028499741721878
$ws.Range("D9").Value = 1.031433594701899
$ws.Range("E9").Value = 1.037633943946112
$ws.Range("F9").Value = 1.046559977135118
$ws.Range("I9").Value = 1.032213562699357
$ws.Range("J9").Value = 1.034346417672459
$ws.Range("K9").Value = 1.034611140992192
$ws.Range("L9").Value = 1.040791200785214
$ws.Range("M9").Value = 1.049688487388761
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.027050557283306
$ws.Range("D10").Value = 1.030412437789532
$ws.Range("E10").Value = 1.036339322357501
$ws.Range("F10").Value = 1.045112469425325
$ws.Range("I10").Value = 1.031991700558728
$ws.Range("J10").Value = 1.033444795627058
$ws.Range("K10").Value = 1.03387811957759
$ws.Range("L10").Value = 1.039783709011029
$ws.Range("M10").Value = 1.04852582306459
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026423759925416
$ws.Range("D11").Value = 1.029970609555492
$ws.Range("E11").Value = 1.035779849378668
$ws.Range("F11").Value = 1.044486912616342
$ws.Range("I11").Value = 1.031894157753258
$ws.Range("J11").Value = 1.033054355811293
$ws.Range("K11").Value = 1.033560209932741
$ws.Range("L11").Value = 1.039347769260921
$ws.Range("M11").Value = 1.048022836310044
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026191046999591
$ws.Range("D12").Value = 1.029806547368405
$ws.Range("E12").Value = 1.035572203891665
$ws.Range("F12").Value = 1.044254738237868
$ws.Range("I12").Value = 1.031857704943848
$ws.Range("J12").Value = 1.032909325355927
$ws.Range("K12").Value = 1.033442049081932
$ws.Range("L12").Value = 1.039185889651985
$ws.Range("M12").Value = 1.047836074184537
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026240959830612
$ws.Range("D13").Value = 1.029841736865153
$ws.Range("E13").Value = 1.035616736942511
$ws.Range("F13").Value = 1.044304532022972
$ws.Range("I13").Value = 1.031865534190353
$ws.Range("J13").Value = 1.032940435033152
$ws.Range("K13").Value = 1.033467398365061
$ws.Range("L13").Value = 1.039220611195386
$ws.Range("M13").Value = 1.047876132126932
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02640452160322
$ws.Range("D14").Value = 1.02995704704449
$ws.Range("E14").Value = 1.035762681908959
$ws.Range("F14").Value = 1.044467717221912
$ws.Range("I14").Value = 1.031891149059266
$ws.Range("J14").Value = 1.03304236761871
$ws.Range("K14").Value = 1.033550444244969
$ws.Range("L14").Value = 1.039334387259834
$ws.Range("M14").Value = 1.048007397068404
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026505311686885
$ws.Range("D15").Value = 1.030028100451952
$ws.Range("E15").Value = 1.035852625682242
$ws.Range("F15").Value = 1.044568285602512
$ws.Range("I15").Value = 1.03190690194784
$ws.Range("J15").Value = 1.033105171185774
$ws.Range("K15").Value = 1.033601601638418
$ws.Range("L15").Value = 1.039404494819095
$ws.Range("M15").Value = 1.048088282984884
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.027092170783403
$ws.Range("D16").Value = 1.030441767753578
$ws.Range("E16").Value = 1.036376476174609
$ws.Range("F16").Value = 1.045154011473182
$ws.Range("I16").Value = 1.031998143122166
$ws.Range("J16").Value = 1.033470707254016
$ws.Range("K16").Value = 1.033899207610997
$ws.Range("L16").Value = 1.039812647541466
$ws.Range("M16").Value = 1.048559214323436
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027460482525935
$ws.Range("D17").Value = 1.030701342257482
$ws.Range("E17").Value = 1.036705371046843
$ws.Range("F17").Value = 1.045521750443442
$ws.Range("I17").Value = 1.032054981771676
$ws.Range("J17").Value = 1.033699990647539
$ws.Range("K17").Value = 1.034085753100488
$ws.Range("L17").Value = 1.040068754974444
$ws.Range("M17").Value = 1.04885473969818
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027675380881922
$ws.Range("D18").Value = 1.030852780325689
$ws.Range("E18").Value = 1.036897316547985
$ws.Range("F18").Value = 1.045736364324057
$ws.Range("I18").Value = 1.032087992438651
$ws.Range("J18").Value = 1.033833724662619
$ws.Range("K18").Value = 1.034194512984862
$ws.Range("L18").Value = 1.040218168022834
$ws.Range("M18").Value = 1.049027158424445
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027748667230739
$ws.Range("D19").Value = 1.03090442228628
$ws.Range("E19").Value = 1.036962783056365
$ws.Range("F19").Value = 1.045809562075544
$ws.Range("I19").Value = 1.032099224060653
$ws.Range("J19").Value = 1.033879323947618
$ws.Range("K19").Value = 1.03423158896129
$ws.Range("L19").Value = 1.040269119051178
$ws.Range("M19").Value = 1.04908595614121
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027420959084891
$ws.Range("D20").Value = 1.030673488981685
$ws.Range("E20").Value = 1.036670072694338
$ws.Range("F20").Value = 1.045482283320484
$ws.Range("I20").Value = 1.032048898243118
$ws.Range("J20").Value = 1.033675391023605
$ws.Range("K20").Value = 1.0340657435898
$ws.Range("L20").Value = 1.040041273968201
$ws.Range("M20").Value = 1.048823028085343
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026356353761758
$ws.Range("D21").Value = 1.029923089609702
$ws.Range("E21").Value = 1.035719700113535
$ws.Range("F21").Value = 1.044419658122051
$ws.Range("I21").Value = 1.03188361221704
$ws.Range("J21").Value = 1.033012351107518
$ws.Range("K21").Value = 1.03352599136664
$ws.Range("L21").Value = 1.039300881715419
$ws.Range("M21").Value = 1.047968740884796
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025687616719573
$ws.Range("D22").Value = 1.029451588178378
$ws.Range("E22").Value = 1.035123134044033
$ws.Range("F22").Value = 1.043752616134901
$ws.Range("I22").Value = 1.031778411566573
$ws.Range("J22").Value = 1.032595450555473
$ws.Range("K22").Value = 1.033186194001796
$ws.Range("L22").Value = 1.038835644975304
$ws.Range("M22").Value = 1.047432019637011
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026042067532747
$ws.Range("D23").Value = 1.029701510652981
$ws.Range("E23").Value = 1.035439292502214
$ws.Range("F23").Value = 1.044106125555524
$ws.Range("I23").Value = 1.031834301453783
$ws.Range("J23").Value = 1.032816459082273
$ws.Range("K23").Value = 1.033366367775896
$ws.Range("L23").Value = 1.039082249117333
$ws.Range("M23").Value = 1.04771650710371
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027438817815514
$ws.Range("D24").Value = 1.030686074577367
$ws.Range("E24").Value = 1.036686022170941
$ws.Range("F24").Value = 1.045500116450171
$ws.Range("I24").Value = 1.03205164756772
$ws.Range("J24").Value = 1.033686506544485
$ws.Range("K24").Value = 1.034074785177714
$ws.Range("L24").Value = 1.040053691359693
$ws.Range("M24").Value = 1.048837357063038
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.02906234632689
$ws.Range("D25").Value = 1.031829881640143
$ws.Range("E25").Value = 1.038136962517277
$ws.Range("F25").Value = 1.047122387426001
$ws.Range("I25").Value = 1.032298279348333
$ws.Range("J25").Value = 1.034696022549613
$ws.Range("K25").Value = 1.034894940373036
$ws.Range("L25").Value = 1.041182165458554
$ws.Range("M25").Value = 1.050139753700381
